$p = $ppt.ActivePresentation

# Update the "Created with GroupDocs.Assembly ..." run text on slide 1.
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$origHeight = $sh.Height

$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2)
$run = $para.Runs(1)
$run.Text = "Created with GroupDocs.Assembly 25.12."

# Editing the run text re-triggers the shape's autofit, which grows its
# height; restore the original (autofit) height so the shape geometry is
# left untouched, matching the source edit.
$sh.Height = $origHeight

# Update the presentation-level Aspose.Slides generator tags.
$p.Tags.Add("AS_NET", "4.0.30319.42000")
$p.Tags.Add("AS_OS", "Microsoft Windows NT 6.2.9200.0")
$p.Tags.Add("AS_TITLE", "Aspose.Slides for .NET 4.0 Client Profile")
